$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 460, shifting existing rows 460:481 down to 461:482
$ws.Rows("460:460").Insert()

# Populate the new row 460 with the new record
$ws.Range("A460").Value = 4
$ws.Range("B460").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C460").Value = "Los Lagos"
$ws.Range("D460").Value = 45147
$ws.Range("E460").Value = 10
$ws.Range("F460").Value = 100112043
$ws.Range("G460").Value = "Pepino ensalada"
$ws.Range("H460").Value = "Sin especificar"
$ws.Range("I460").Value = "Primera"
$ws.Range("J460").Value = 80
$ws.Range("K460").Value = 15000
$ws.Range("L460").Value = 15000
$ws.Range("M460").Value = 15000
$ws.Range("N460").Value = '$/caja 60 unidades'
$ws.Range("O460").Value = "Región de Arica y Parinacota"
$ws.Range("P460").Value = 250
$ws.Range("Q460").Value = 60
$ws.Range("R460").Value = "Hortaliza"
